$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry updates a single cell's displayed/stored text. The source data
# (coinranking.com price/volume snapshot) is refreshed by a scheduled
# GitHub Actions job; cells are plain text in the original workbook (many
# "prices" use '.' as a thousands separator rather than a decimal point, so
# they are not real numbers), hence every value below is entered with a
# leading apostrophe to force text entry and ClearFormats() is called right
# after so Excel does not leave a lingering "Text" / quote-prefix style on
# the cell (keeping the cell style identical to the untouched original).
$updates = @(
    @{ Cell = "D2"; Value = "'26.261.16" }
    @{ Cell = "E2"; Value = "'  -0.07%  " }
    @{ Cell = "D3"; Value = "'1.592.74" }
    @{ Cell = "E3"; Value = "'  +0.10%  " }
    @{ Cell = "E4"; Value = "'  -0.01%  " }
    @{ Cell = "D5"; Value = "'212.73" }
    @{ Cell = "E5"; Value = "'  -0.17%  " }
    @{ Cell = "E6"; Value = "'  -0.18%  " }
    @{ Cell = "E7"; Value = "'  +0.03%  " }
    @{ Cell = "E8"; Value = "'  -0.53%  " }
    @{ Cell = "E9"; Value = "'  -0.44%  " }
    @{ Cell = "D10"; Value = "'18.95" }
    @{ Cell = "E10"; Value = "'  -2.14%  " }
    @{ Cell = "D11"; Value = "'0.0851" }
    @{ Cell = "E11"; Value = "'  +0.20%  " }
    @{ Cell = "E12"; Value = "'  +0.08%  " }
    @{ Cell = "D13"; Value = "'1.613.35" }
    @{ Cell = "E13"; Value = "'  +1.42%  " }
    @{ Cell = "E14"; Value = "'  -1.08%  " }
    @{ Cell = "E15"; Value = "'  -2.71%  " }
    @{ Cell = "D16"; Value = "'63.93" }
    @{ Cell = "E16"; Value = "'  -0.86%  " }
    @{ Cell = "D17"; Value = "'26.260.08" }
    @{ Cell = "E17"; Value = "'  -0.07%  " }
    @{ Cell = "D18"; Value = "'0.0₃0723" }
    @{ Cell = "E18"; Value = "'  -0.59%  " }
    @{ Cell = "D19"; Value = "'215.26" }
    @{ Cell = "E19"; Value = "'  +0.96%  " }
    @{ Cell = "D20"; Value = "'7.39" }
    @{ Cell = "E20"; Value = "'  -1.22%  " }
    @{ Cell = "E21"; Value = "'  +0.10%  " }
    @{ Cell = "E22"; Value = "'  -0.09%  " }
    @{ Cell = "E23"; Value = "'  -0.01%  " }
    @{ Cell = "E24"; Value = "'  -2.59%  " }
    @{ Cell = "D25"; Value = "'144.78" }
    @{ Cell = "E25"; Value = "'  -0.07%  " }
    @{ Cell = "E26"; Value = "'  +0.03%  " }
    @{ Cell = "E27"; Value = "'  -1.36%  " }
    @{ Cell = "E28"; Value = "'  +0.63%  " }
    @{ Cell = "E29"; Value = "'  -0.52%  " }
    @{ Cell = "D30"; Value = "'0.0495" }
    @{ Cell = "E30"; Value = "'  -1.04%  " }
    @{ Cell = "E31"; Value = "'  -0.01%  " }
    @{ Cell = "E32"; Value = "'  -0.45%  " }
    @{ Cell = "D33"; Value = "'1.426.97" }
    @{ Cell = "E33"; Value = "'  +6.60%  " }
    @{ Cell = "E34"; Value = "'  +0.01%  " }
    @{ Cell = "E35"; Value = "'  -0.97%  " }
    @{ Cell = "E36"; Value = "'  -1.18%  " }
    @{ Cell = "D37"; Value = "'0.566" }
    @{ Cell = "E37"; Value = "'  -4.56%  " }
    @{ Cell = "E38"; Value = "'  -0.54%  " }
    @{ Cell = "E39"; Value = "'  +0.68%  " }
    @{ Cell = "D40"; Value = "'5.76" }
    @{ Cell = "E40"; Value = "'  -0.26%  " }
    @{ Cell = "E41"; Value = "'  +0.04%  " }
    @{ Cell = "E42"; Value = "'  +0.87%  " }
    @{ Cell = "D43"; Value = "'0.922" }
    @{ Cell = "E43"; Value = "'  -8.36%  " }
    @{ Cell = "D44"; Value = "'0.760" }
    @{ Cell = "E44"; Value = "'  -0.35%  " }
    @{ Cell = "D45"; Value = "'1.729.42" }
    @{ Cell = "E45"; Value = "'  +0.24%  " }
    @{ Cell = "D46"; Value = "'60.81" }
    @{ Cell = "E46"; Value = "'  -1.72%  " }
    @{ Cell = "D47"; Value = "'86.67" }
    @{ Cell = "E47"; Value = "'  +0.17%  " }
    @{ Cell = "E48"; Value = "'  +0.47%  " }
    @{ Cell = "E49"; Value = "'  -1.60%  " }
    @{ Cell = "E50"; Value = "'  -0.79%  " }
    @{ Cell = "E51"; Value = "'  -3.14%  " }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    $cell.Value = $u.Value
    $cell.ClearFormats()
}
